$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.571.06'
$ws.Range("E2").Value = '  -2.33%  '
$ws.Range("D3").Value = '2.892.59'
$ws.Range("E3").Value = '  -2.04%  '
$ws.Range("E4").Value = '  -0.02%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '567.86'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -4.22%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '143.03'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -3.60%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("D9").Value = '2.890.95'
$ws.Range("E9").Value = '  -2.02%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '6.92'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.44%  '
$ws.Range("E11").Value = '  -3.12%  '
$ws.Range("E12").Value = '  -2.47%  '
$ws.Range("E13").Value = '  -1.86%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '31.77'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.13%  '
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").Value = '3.372.07'
$ws.Range("E16").Value = '  -2.05%  '
$ws.Range("D17").Value = '61.561.58'
$ws.Range("E17").Value = '  -2.36%  '
$ws.Range("E18").Value = '  -2.28%  '
$ws.Range("D19").Value = '2.889.55'
$ws.Range("E19").Value = '  -2.18%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '431.71'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.14%  '
$ws.Range("E21").Value = '  -3.44%  '
$ws.Range("E22").Value = '  -2.26%  '
$ws.Range("E23").Value = '  -2.69%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '79.17'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -2.04%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '11.87'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -12.38%  '
$ws.Range("E28").Value = '  -6.26%  '
$ws.Range("E29").Value = '  +2.88%  '
$ws.Range("E30").Value = '  -4.49%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '2.49'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -4.33%  '
$ws.Range("E32").Value = '  -9.68%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  -2.14%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '25.49'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -3.58%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.956'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -3.57%  '
$ws.Range("E37").Value = '  -4.35%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '48.89'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.78%  '
$ws.Range("E39").Value = '  -5.36%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.80'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -9.58%  '
$ws.Range("E41").Value = '  -3.49%  '
$ws.Range("E42").Value = '  -3.44%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '39.63'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("E44").Value = '  -4.66%  '
$ws.Range("D45").Value = '2.688.22'
$ws.Range("E45").Value = '  -0.77%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '132.76'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -2.28%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.0334'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.23%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '343.46'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -4.93%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("E51").Value = '  -5.57%  '
